$wb = $excel.ActiveWorkbook

# Map of sheet -> row -> {column letter -> new value ($null means clear/blank the cell)}
# Applies the numeric corrections from the scheduled-runner price refresh.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2119.9473
$ws.Range("J17").Value = 2119.9473
$ws.Range("L17").Value = 6359.841899999999
$ws.Range("N17").Value = -6695.841899999999

$ws.Range("H32").Value = 2891.111
$ws.Range("I32").Value = 2147.25
$ws.Range("J32").Value = 3486.2
$ws.Range("K32").Value = 2147.25
$ws.Range("L32").Value = 3486.2
$ws.Range("M32").Value = -1821.25
$ws.Range("N32").Value = -4138.2

$ws.Range("H70").Value = 3551.7568
$ws.Range("I70").Value = 1636.5385
$ws.Range("J70").Value = 4589.1665
$ws.Range("K70").Value = 4909.6155
$ws.Range("L70").Value = 13767.4995
$ws.Range("M70").Value = -4639.6155
$ws.Range("N70").Value = -14307.4995

$ws.Range("H73").Value = 3551.7568
$ws.Range("I73").Value = 1636.5385
$ws.Range("J73").Value = 4589.1665
$ws.Range("K73").Value = 4909.6155
$ws.Range("L73").Value = 13767.4995
$ws.Range("M73").Value = -3973.6155
$ws.Range("N73").Value = -15639.4995

$ws.Range("H76").Value = 20003300
$ws.Range("I76").Value = 33335966
$ws.Range("K76").Value = 33335966
$ws.Range("M76").Value = -33335651

$ws.Range("H79").Value = 20003300
$ws.Range("I79").Value = 33335966
$ws.Range("K79").Value = 33335966
$ws.Range("M79").Value = -33334874

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").Value = $null

$ws.Range("H116").Value = 4000
$ws.Range("I116").Value = 4000
$ws.Range("K116").Value = 4000
$ws.Range("M116").Value = -558

$ws.Range("H132").Value = 6143.727
$ws.Range("I132").Value = 6559.1
$ws.Range("J132").Value = 1990
$ws.Range("K132").Value = 19677.3
$ws.Range("L132").Value = 5970
$ws.Range("M132").Value = -17147.3
$ws.Range("N132").Value = -11030

$ws.Range("H138").Value = 790.7059
$ws.Range("J138").Value = 1425
$ws.Range("L138").Value = 4275
$ws.Range("N138").Value = -14555


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 4851.3335
$ws.Range("I110").Value = 5201.4546
$ws.Range("K110").Value = 5201.4546
$ws.Range("M110").Value = -3156.4546


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 2202
$ws.Range("J64").Value = 2563.8333
$ws.Range("L64").Value = 2563.8333
$ws.Range("N64").Value = -3013.8333

$ws.Range("H67").Value = 2202
$ws.Range("J67").Value = 2563.8333
$ws.Range("L67").Value = 2563.8333
$ws.Range("N67").Value = -4123.8333

$ws.Range("H86").Value = 2353.4546
$ws.Range("I86").Value = 2465.4443
$ws.Range("K86").Value = 2465.4443
$ws.Range("M86").Value = -1342.4443

$ws.Range("H88").Value = 37499.5
$ws.Range("J88").Value = 37499.5
$ws.Range("L88").Value = 37499.5
$ws.Range("N88").Value = -38311.5

$ws.Range("H89").Value = 2353.4546
$ws.Range("I89").Value = 2465.4443
$ws.Range("K89").Value = 12327.2215
$ws.Range("M89").Value = -6711.2215

$ws.Range("H91").Value = 37499.5
$ws.Range("J91").Value = 37499.5
$ws.Range("L91").Value = 37499.5
$ws.Range("N91").Value = -40307.5


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 8.333333
$ws.Range("I23").Value = 9
$ws.Range("J23").Value = 8
$ws.Range("K23").Value = 9
$ws.Range("L23").Value = 8
$ws.Range("M23").Value = 231
$ws.Range("N23").Value = -488

$ws.Range("H27").Value = 8.333333
$ws.Range("I27").Value = 9
$ws.Range("J27").Value = 8
$ws.Range("K27").Value = 9
$ws.Range("L27").Value = 8
$ws.Range("M27").Value = 183
$ws.Range("N27").Value = -392

$ws.Range("H63").Value = 52567.75
$ws.Range("J63").Value = 75135.5
$ws.Range("L63").Value = 75135.5
$ws.Range("N63").Value = -76507.5

$ws.Range("H66").Value = 52567.75
$ws.Range("J66").Value = 75135.5
$ws.Range("L66").Value = 225406.5
$ws.Range("N66").Value = -232270.5


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 18353.555
$ws.Range("I7").Value = 25567.166
$ws.Range("J7").Value = 3926.3333
$ws.Range("K7").Value = 76701.49800000001
$ws.Range("L7").Value = 11778.9999
$ws.Range("M7").Value = -76589.49800000001
$ws.Range("N7").Value = -12002.9999

$ws.Range("H64").Value = 353.5
$ws.Range("I64").Value = 353
$ws.Range("K64").Value = 1059
$ws.Range("M64").Value = -789

$ws.Range("H67").Value = 353.5
$ws.Range("I67").Value = 353
$ws.Range("K67").Value = 1059
$ws.Range("M67").Value = -123

$ws.Range("H129").Value = 2656.6
$ws.Range("I129").Value = 700
$ws.Range("J129").Value = 3145.75
$ws.Range("K129").Value = 2100
$ws.Range("L129").Value = 9437.25
$ws.Range("M129").Value = 2900
$ws.Range("N129").Value = -19437.25


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 6312.5
$ws.Range("I43").Value = 916.6667
$ws.Range("J43").Value = 22500
$ws.Range("K43").Value = 916.6667
$ws.Range("L43").Value = 22500
$ws.Range("M43").Value = -765.6667
$ws.Range("N43").Value = -22802

$ws.Range("H46").Value = 24999
$ws.Range("J46").Value = 24999
$ws.Range("L46").Value = 24999
$ws.Range("N46").Value = -25311

$ws.Range("H70").Value = 40739.582
$ws.Range("I70").Value = 43248.1
$ws.Range("K70").Value = 43248.1
$ws.Range("M70").Value = -42978.1

$ws.Range("H73").Value = 40739.582
$ws.Range("I73").Value = 43248.1
$ws.Range("K73").Value = 43248.1
$ws.Range("M73").Value = -42312.1

$ws.Range("H80").Value = 7333
$ws.Range("J80").Value = 6499.5
$ws.Range("L80").Value = 6499.5
$ws.Range("N80").Value = -8495.5

$ws.Range("H83").Value = 7333
$ws.Range("J83").Value = 6499.5
$ws.Range("L83").Value = 32497.5
$ws.Range("N83").Value = -42481.5

$ws.Range("H86").Value = 111565
$ws.Range("J86").Value = 143000
$ws.Range("L86").Value = 143000
$ws.Range("N86").Value = -145372

$ws.Range("H89").Value = 111565
$ws.Range("J89").Value = 143000
$ws.Range("L89").Value = 429000
$ws.Range("N89").Value = -440856

$ws.Range("H122").Value = 1406.5
$ws.Range("I122").Value = 1406.5
$ws.Range("K122").Value = 4219.5
$ws.Range("M122").Value = -1769.5


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 983.625
$ws.Range("J46").Value = 600
$ws.Range("L46").Value = 600
$ws.Range("N46").Value = -976

$ws.Range("H68").Value = 3540.8
$ws.Range("I68").Value = 2400.3333
$ws.Range("K68").Value = 2400.3333
$ws.Range("M68").Value = -1651.3333

$ws.Range("H71").Value = 3540.8
$ws.Range("I71").Value = 2400.3333
$ws.Range("K71").Value = 12001.6665
$ws.Range("M71").Value = -8257.666499999999


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12329
$ws.Range("J62").Value = 12329
$ws.Range("L62").Value = 12329
$ws.Range("N62").Value = -13577

$ws.Range("H65").Value = 12329
$ws.Range("J65").Value = 12329
$ws.Range("L65").Value = 61645
$ws.Range("N65").Value = -67885

$ws.Range("H122").Value = 5770
$ws.Range("I122").Value = 6148
$ws.Range("K122").Value = 18444
$ws.Range("M122").Value = -15994

$ws.Range("H132").Value = 4283.2666
$ws.Range("I132").Value = 2479.0833
$ws.Range("K132").Value = 7437.249899999999
$ws.Range("M132").Value = -4907.249899999999

